# TSTool-Command-List.xlsx update:
#  - Add a new "ReadExcelWorkbook" command row, inserted in alphabetical
#    order right before "ReadHecDss" (which was row 106), pushing every
#    row below it down by one.
#  - Mark the "ChangePeriod" command row (row 13) as Y/Y for the two
#    boolean columns (it previously had no value there).
#  - Refresh the view: scroll position and active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row, shifting rows 106:227 down to 107:228 -----------
$ws.Rows.Item(106).Insert()

# New row 106 content: the ReadExcelWorkbook command.
$ws.Range("A106").Value = "ReadExcelWorkbook"
$ws.Range("B106").Value = "Read an Excel workbook into memory so that it can be manipulated."
$ws.Range("C106").Value = "Y"
$ws.Range("D106").Value = "Y"

# Match the centered formatting used by the other Y/NA cells in C/D.
$ws.Range("C106").HorizontalAlignment = -4108
$ws.Range("D106").HorizontalAlignment = -4108

# --- ChangePeriod (row 13) now handles processor properties/For loops ----
$ws.Range("C13").Value = "Y"
$ws.Range("D13").Value = "Y"
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("D13").HorizontalAlignment = -4108

# --- Refresh the saved view/selection state -------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 201
$win.ScrollColumn = 1
$ws.Range("C108").Select()
